$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B5").Value = "Four age groups"
$ws.Range("D5").Value = "[60, 65], (65, 70], (70, 75],  (75, 80]"
$ws.Range("G13").Value = "LBXBPB"
